{"js": "// Append the new sentences to the end of the document body (same paragraph),\n// matching the authored diff which adds runs of text after the existing\n// run but inside the very same <w:p> (no new paragraph is introduced).\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = body.paragraphs.items[body.paragraphs.items.length - 1];\n\nconst addition =\n  \"They were ruled by a ferocious and ever-hungry lion.\" +\n  \" The animals all lived in fear of their king, the lion as he mercilessly hunted the animals to satisfy his hunger.\";\n\nlastParagraph.insertText(addition, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Append the new sentences to the end of the document (same paragraph),\n# matching the authored diff which adds runs of text after the existing\n# run but inside the very same paragraph (no new paragraph is introduced).\n$d = $word.ActiveDocument\n\n$addition = \"They were ruled by a ferocious and ever-hungry lion.\" + \" The animals all lived in fear of their king, the lion as he mercilessly hunted the animals to satisfy his hunger.\"\n\n$para = $d.Paragraphs.Last\n$r = $para.Range\n$r.Collapse($wdCollapseEnd)\n$r.InsertAfter($addition)\n"}
